$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rows 12 & 13 used to describe the separate "Speaker_L" / "Speaker_R" tests.
# They are now combined into a single generic "Speaker" test, which replaces
# row 12 and keeps the same ID/N1 numbering (11).
$ws.Range("B12").Value = "Speaker"
$ws.Range("H12").Value = "Speaker"
$ws.Range("J12").Value = "C:\TestProgram\0WM\0WM_BLT3\Test_Speaker.bat"
$ws.Range("K12").Value = "C:\TestProgram\0WM\0WM_BLT3\log\Test_Speaker_CheckLog.bat"

# The now-redundant "Speaker_R" row (13) is removed; this pulls the old
# "LED10to16" row (14) up to become the new row 13, unchanged otherwise.
$ws.Rows("13").Delete()

# Tidy up the view state to match the saved workbook.
$ws.Range("K19").Select()
